$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 46
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '11776.00'

# Row 9
$ws.Range("C9").Value = 82
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '38704.00'

# Row 10
$ws.Range("A10").Value = ''
$ws.Range("C10").Value = 45
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.0'
$ws.Range("E10").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F10").Value = 0
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '0.00'

# Row 11
$ws.Range("A11").Value = 'P. point'
$ws.Range("C11").Value = 84
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6'
$ws.Range("E11").Value = 'On board'
$ws.Range("F11").Value = 136
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '11424.00'

# Row 12
$ws.Range("C12").Value = 57
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.0'
$ws.Range("E12").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 23
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '1311.00'

# Row 13
$ws.Range("C13").Value = 11
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.0'
$ws.Range("E13").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 50
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '550.00'

# Row 14
$ws.Range("C14").Value = 13
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.0'
$ws.Range("E14").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 78
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '1014.00'

# Row 15
$ws.Range("A15").Value = 'Each'
$ws.Range("C15").Value = 47
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.0'
$ws.Range("E15").Value = 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = 30
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '1410.00'

# Row 16
$ws.Range("A16").Value = 'Each'
$ws.Range("C16").Value = 51
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '9.0'
$ws.Range("E16").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 219
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '11169.00'

# Row 17
$ws.Range("A17").Value = 'Each'
$ws.Range("C17").Value = 13
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '10.0'
$ws.Range("E17").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F17").Value = 303
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '3939.00'

# Row 18
$ws.Range("A18").Value = ''
$ws.Range("C18").Value = 50
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.0'
$ws.Range("E18").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F18").Value = 0
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '0.00'

# Row 19
$ws.Range("A19").Value = 'R. mtr.'
$ws.Range("C19").Value = 57
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17'
$ws.Range("E19").Value = '25 mm'
$ws.Range("F19").Value = 56
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '3192.00'

# Row 20
$ws.Range("C20").Value = 74
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19'
$ws.Range("E20").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F20").Value = 81
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '5994.00'

# Row 21
$ws.Range("A21").Value = 'Mtr.'
$ws.Range("C21").Value = 10
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '20'
$ws.Range("E21").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F21").Value = 122
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '1220.00'

# Row 22
$ws.Range("A22").Value = 'Set'
$ws.Range("C22").Value = 98
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.0'
$ws.Range("E22").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F22").Value = 5733
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '561834.00'

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.0'
$ws.Range("E23").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '

# Row 24
$ws.Range("A24").Value = 'Mtr.'
$ws.Range("C24").Value = 98
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23'
$ws.Range("E24").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F24").Value = 20
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '1960.00'

# Row 25
$ws.Range("A25").Value = ''
$ws.Range("C25").Value = 63
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.0'
$ws.Range("E25").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F25").Value = 0
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '0.00'

# Row 26
$ws.Range("C26").Value = 42
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '25'
$ws.Range("E26").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F26").Value = 1890
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '79380.00'

# Row 27
$ws.Range("A27").Value = 'Each'
$ws.Range("C27").Value = 62
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27'
$ws.Range("E27").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F27").Value = 492
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '30504.00'

# Row 28
$ws.Range("C28").Value = 29
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.0'
$ws.Range("E28").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'

# Row 29
$ws.Range("C29").Value = 31
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34'
$ws.Range("E29").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'

# Row 30
$ws.Range("C30").Value = 44

# Row 31
$ws.Range("C31").Value = 4

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '765381.00'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '765381.00'

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '765381.00'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '765381.00'
